$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "36.673.01"
$ws.Range("E2").Value = "  -1.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.008.93"
$ws.Range("E3").Value = "  -0.36%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.00%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.600"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.02%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.76"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.10%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -3.72%  "

# Row 10 - OKB
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.13%  "

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.97%  "

# Row 12 - TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.102"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.301.27"
$ws.Range("E13").Value = "  -0.49%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.66%  "

# Row 15 - Avalanche
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.91%  "

# Row 16 - Polygon
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.755"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.59%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -2.96%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.004.74"
$ws.Range("E18").Value = "  -0.68%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "36.567.56"
$ws.Range("E19").Value = "  -2.02%  "

# Row 20 - Litecoin
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.42%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0796"
$ws.Range("E21").Value = "  -5.18%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.45%  "

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.22%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.09%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +2.37%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  -8.22%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.20%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  -3.95%  "

# Row 29 - ImmutableX
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.33%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  -1.73%  "

# Row 31 - EthereumClassic
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "18.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.41%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -5.70%  "

# Row 34 - Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.39%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.38%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.58%  "

# Row 37 - BinanceUSD
$ws.Range("E37").Value = "  +0.07%  "

# Row 38 - RenderToken -> WEMIXToken
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.47%  "

# Row 39 - WEMIXToken -> RenderToken
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.07%  "

# Row 40 - THORChain
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.68%  "

# Row 41 - HuobiToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.30%  "

# Row 42 - Maker
$ws.Range("D42").Value = "1.457.59"
$ws.Range("E42").Value = "  +1.64%  "

# Row 43 - Cronos
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0929"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.29%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -4.08%  "

# Row 45 - Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "

# Row 46 - TrustWalletToken
$ws.Range("E46").Value = "  -8.17%  "

# Row 47 - InjectiveProtocol
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.38%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  -3.47%  "

# Row 49 - FTXToken
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +26.83%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  -1.73%  "

# Row 51 - FraxShare
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.80%  "

